$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 16.33975766666667
$ws.Range("H2").Value = 49.01927300000001
$ws.Range("I2").Value = 0.2979519994155143
$ws.Range("J2").Value = 0.2979519994155143
$ws.Range("M2").Value = 61.156892
$ws.Range("N2").Value = 183.470676
$ws.Range("O2").Value = 0.9308124486389074
$ws.Range("P2").Value = 0.9308124486389074
$ws.Range("Q2").Value = 999.2887949265054
$ws.Range("R2").Value = 8993.599154338548
$ws.Range("S2").Value = 0.2773374301528132
$ws.Range("T2").Value = 0.2773374301528132
$ws.Range("G3").Value = 16.33975766666667
$ws.Range("H3").Value = 49.01927300000001
$ws.Range("I3").Value = 0.2979519994155143
$ws.Range("J3").Value = 0.2979519994155143
$ws.Range("O3").Value = 0.02171808228502914
$ws.Range("P3").Value = 0.02171808228502914
$ws.Range("Q3").Value = 23.31579933901444
$ws.Range("R3").Value = 209.84219405113
$ws.Range("S3").Value = 0.006470946040295095
$ws.Range("T3").Value = 0.006470946040295094
$ws.Range("G4").Value = 16.33975766666667
$ws.Range("H4").Value = 49.01927300000001
$ws.Range("I4").Value = 0.2979519994155143
$ws.Range("J4").Value = 0.2979519994155143
$ws.Range("M4").Value = 2.00294
$ws.Range("N4").Value = 6.00882
$ws.Range("O4").Value = 0.03048489589491914
$ws.Range("P4").Value = 0.03048489589491914
$ws.Range("Q4").Value = 32.72755422087334
$ws.Range("R4").Value = 294.54798798786
$ws.Range("S4").Value = 0.009083035683864963
$ws.Range("T4").Value = 0.009083035683864963
$ws.Range("G5").Value = 16.33975766666667
$ws.Range("H5").Value = 49.01927300000001
$ws.Range("I5").Value = 0.2979519994155143
$ws.Range("J5").Value = 0.2979519994155143
$ws.Range("M5").Value = 1.115932333333334
$ws.Range("N5").Value = 3.347797
$ws.Range("O5").Value = 0.01698457318114416
$ws.Range("P5").Value = 0.01698457318114415
$ws.Range("Q5").Value = 18.23406389906456
$ws.Range("R5").Value = 164.106575091581
$ws.Range("S5").Value = 0.005060587538541024
$ws.Range("T5").Value = 0.005060587538541023
$ws.Range("I6").Value = 0.2656466977818992
$ws.Range("J6").Value = 0.2656466977818992
$ws.Range("M6").Value = 61.156892
$ws.Range("N6").Value = 183.470676
$ws.Range("O6").Value = 0.9308124486389074
$ws.Range("P6").Value = 0.9308124486389074
$ws.Range("Q6").Value = 890.9413899669145
$ws.Range("R6").Value = 8018.472509702231
$ws.Range("S6").Value = 0.2472672532352094
$ws.Range("T6").Value = 0.2472672532352094
$ws.Range("I7").Value = 0.2656466977818992
$ws.Range("J7").Value = 0.2656466977818992
$ws.Range("O7").Value = 0.02171808228502914
$ws.Range("P7").Value = 0.02171808228502914
$ws.Range("S7").Value = 0.005769336841173555
$ws.Range("T7").Value = 0.005769336841173554
$ws.Range("I8").Value = 0.2656466977818992
$ws.Range("J8").Value = 0.2656466977818992
$ws.Range("M8").Value = 2.00294
$ws.Range("N8").Value = 6.00882
$ws.Range("O8").Value = 0.03048489589491914
$ws.Range("P8").Value = 0.03048489589491914
$ws.Range("Q8").Value = 29.17908496102666
$ws.Range("R8").Value = 262.61176464924
$ws.Range("S8").Value = 0.008098211926710246
$ws.Range("T8").Value = 0.008098211926710244
$ws.Range("I9").Value = 0.2656466977818992
$ws.Range("J9").Value = 0.2656466977818992
$ws.Range("M9").Value = 1.115932333333334
$ws.Range("N9").Value = 3.347797
$ws.Range("O9").Value = 0.01698457318114416
$ws.Range("P9").Value = 0.01698457318114415
$ws.Range("Q9").Value = 16.25704432738378
$ws.Range("R9").Value = 146.313398946454
$ws.Range("S9").Value = 0.004511895778805952
$ws.Range("T9").Value = 0.004511895778805951
$ws.Range("G10").Value = 22.15292366666667
$ws.Range("H10").Value = 66.458771
$ws.Range("I10").Value = 0.4039538427701242
$ws.Range("J10").Value = 0.4039538427701242
$ws.Range("M10").Value = 61.156892
$ws.Range("N10").Value = 183.470676
$ws.Range("O10").Value = 0.9308124486389074
$ws.Range("P10").Value = 0.9308124486389074
$ws.Range("Q10").Value = 1354.803960166577
$ws.Range("R10").Value = 12193.2356414992
$ws.Range("S10").Value = 0.3760052655259556
$ws.Range("T10").Value = 0.3760052655259556
$ws.Range("G11").Value = 22.15292366666667
$ws.Range("H11").Value = 66.458771
$ws.Range("I11").Value = 0.4039538427701242
$ws.Range("J11").Value = 0.4039538427701242
$ws.Range("O11").Value = 0.02171808228502914
$ws.Range("P11").Value = 0.02171808228502914
$ws.Range("Q11").Value = 31.61081905383444
$ws.Range("R11").Value = 284.49737148451
$ws.Range("S11").Value = 0.008773102796635283
$ws.Range("T11").Value = 0.008773102796635282
$ws.Range("G12").Value = 22.15292366666667
$ws.Range("H12").Value = 66.458771
$ws.Range("I12").Value = 0.4039538427701242
$ws.Range("J12").Value = 0.4039538427701242
$ws.Range("M12").Value = 2.00294
$ws.Range("N12").Value = 6.00882
$ws.Range("O12").Value = 0.03048489589491914
$ws.Range("P12").Value = 0.03048489589491914
$ws.Range("Q12").Value = 44.37097692891334
$ws.Range("R12").Value = 399.33879236022
$ws.Range("S12").Value = 0.01231449084319977
$ws.Range("T12").Value = 0.01231449084319977
$ws.Range("G13").Value = 22.15292366666667
$ws.Range("H13").Value = 66.458771
$ws.Range("I13").Value = 0.4039538427701242
$ws.Range("J13").Value = 0.4039538427701242
$ws.Range("M13").Value = 1.115932333333334
$ws.Range("N13").Value = 3.347797
$ws.Range("O13").Value = 0.01698457318114416
$ws.Range("P13").Value = 0.01698457318114415
$ws.Range("Q13").Value = 24.72116379749856
$ws.Range("R13").Value = 222.490474177487
$ws.Range("S13").Value = 0.006860983604333576
$ws.Range("T13").Value = 0.006860983604333575
$ws.Range("G14").Value = 1.779426333333333
$ws.Range("H14").Value = 5.338279
$ws.Range("I14").Value = 0.03244746003246218
$ws.Range("J14").Value = 0.03244746003246217
$ws.Range("M14").Value = 61.156892
$ws.Range("N14").Value = 183.470676
$ws.Range("O14").Value = 0.9308124486389074
$ws.Range("P14").Value = 0.9308124486389074
$ws.Range("Q14").Value = 108.8241840896227
$ws.Range("R14").Value = 979.417656806604
$ws.Range("S14").Value = 0.0302024997249292
$ws.Range("T14").Value = 0.03020249972492919
$ws.Range("G15").Value = 1.779426333333333
$ws.Range("H15").Value = 5.338279
$ws.Range("I15").Value = 0.03244746003246218
$ws.Range("J15").Value = 0.03244746003246217
$ws.Range("O15").Value = 0.02171808228502914
$ws.Range("P15").Value = 0.02171808228502914
$ws.Range("Q15").Value = 2.539128680665555
$ws.Range("R15").Value = 22.85215812599
$ws.Range("S15").Value = 0.0007046966069252079
$ws.Range("T15").Value = 0.0007046966069252077
$ws.Range("G16").Value = 1.779426333333333
$ws.Range("H16").Value = 5.338279
$ws.Range("I16").Value = 0.03244746003246218
$ws.Range("J16").Value = 0.03244746003246217
$ws.Range("M16").Value = 2.00294
$ws.Range("N16").Value = 6.00882
$ws.Range("O16").Value = 0.03048489589491914
$ws.Range("P16").Value = 0.03048489589491914
$ws.Range("Q16").Value = 3.564084180086667
$ws.Range("R16").Value = 32.07675762078
$ws.Range("S16").Value = 0.0009891574411441592
$ws.Range("T16").Value = 0.0009891574411441588
$ws.Range("G17").Value = 1.779426333333333
$ws.Range("H17").Value = 5.338279
$ws.Range("I17").Value = 0.03244746003246218
$ws.Range("J17").Value = 0.03244746003246217
$ws.Range("M17").Value = 1.115932333333334
$ws.Range("N17").Value = 3.347797
$ws.Range("O17").Value = 0.01698457318114416
$ws.Range("P17").Value = 0.01698457318114415
$ws.Range("Q17").Value = 1.985719380151445
$ws.Range("R17").Value = 17.871474421363
$ws.Range("S17").Value = 0.000551106259463604
$ws.Range("T17").Value = 0.0005511062594636038
